# Loan RBI, Variable Instalments
# The "Repayment Schedule" sheet gains a new (blank) column, inserted
# immediately before the existing "Late" column, so a "Variable Instalments"
# field can be tracked going forward. Everything from that column onward
# (Late / Due / Outstanding) shifts one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at N (pushes Late/Due/Outstanding -> O/P/Q).
$ws.Range("N1").EntireColumn.Insert()

# New column keeps a plain (non-bestFit) width of 10 characters.
$ws.Columns("N").ColumnWidth = 9.14

# Update the remembered selection on the sheet.
$ws.Range("T8").Select()
